$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "Growth and yield of 'Muscat Hamburg' grape was compared when grafted on 'Dog Ridge' rootstock and self-rooted cuttings in Tamil Nadu, India. Grafted vines showed better results."

$ws.Range("D6").Value = "Silicon preparations can promote growth in forest seedlings, increasing biomass and improving nutritional value. Spraying with a 2% concentration is most effective, particularly for oak seedlings."

$ws.Range("D7").Value = "This study investigated the effects of different soil media mixtures containing phosphogypsum on the growth of young pine seedlings. The mixtures did not have harmful effects, but longer-term observations are needed. Lower dosages appear to be more promising and cost-effective. Further testing is recommended for heavy metals and microbiome changes."

$ws.Range("D9").Value = "This study investigates the growth-survival trade-offs in non-phanerophyte species used in a coastal dune restoration project. The results suggest that plant species of foredune communities have higher growth but lower survival rates, providing insights for cost-effective ecosystem restoration actions."

$ws.Range("D10").Value = "A study found that the colonization of fauna in eelgrass restoration plots was rapid and similar regardless of patch size, suggesting smaller patches can be as effective for promoting biodiversity."

$ws.Range("C11").Value = "No"
$ws.Range("D11").Value = ""
